$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells that would otherwise be auto-converted to numbers by Excel,
# then restore default (Normal) style so no stray style index is introduced.
function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '73.302.09'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '3.976.58'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue $ws 'D5' '616.75'
$ws.Range('E5').Value = '  +14.19%  '
Set-TextValue $ws 'D6' '166.65'
$ws.Range('E6').Value = '  +9.79%  '
Set-TextValue $ws 'D7' '0.681'
$ws.Range('E7').Value = '  -2.45%  '
$ws.Range('E8').Value = '  -0.03%  '
Set-TextValue $ws 'D9' '0.757'
$ws.Range('E9').Value = '  +0.42%  '
Set-TextValue $ws 'D10' '0.186'
$ws.Range('E10').Value = '  +8.17%  '
Set-TextValue $ws 'D11' '55.55'
$ws.Range('E11').Value = '  +2.95%  '
$ws.Range('E12').Value = '  +1.86%  '
Set-TextValue $ws 'D13' '11.11'
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('D14').Value = '4.610.03'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').Value = '3.970.60'
$ws.Range('E15').Value = '  -1.67%  '
$ws.Range('E16').Value = '  +3.35%  '
Set-TextValue $ws 'D17' '14.00'
$ws.Range('E17').Value = '  -2.34%  '
Set-TextValue $ws 'D18' '20.46'
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').Value = '73.055.27'
$ws.Range('E19').Value = '  +1.24%  '
$ws.Range('E20').Value = '  -0.61%  '
Set-TextValue $ws 'D21' '440.34'
$ws.Range('E21').Value = '  -2.25%  '
Set-TextValue $ws 'D22' '4.86'
$ws.Range('E22').Value = '  +14.09%  '
Set-TextValue $ws 'D23' '96.02'
$ws.Range('E23').Value = '  -1.81%  '
$ws.Range('E24').Value = '  -4.18%  '
Set-TextValue $ws 'D25' '14.15'
$ws.Range('E25').Value = '  -3.01%  '
Set-TextValue $ws 'D26' '4.05'
$ws.Range('E26').Value = '  -5.58%  '
Set-TextValue $ws 'D27' '11.02'
$ws.Range('E27').Value = '  -2.14%  '
Set-TextValue $ws 'D28' '5.96'
$ws.Range('E28').Value = '  +0.06%  '
Set-TextValue $ws 'D29' '10.53'
$ws.Range('E29').Value = '  -2.52%  '
Set-TextValue $ws 'D30' '36.14'
$ws.Range('E30').Value = '  -2.64%  '
Set-TextValue $ws 'D31' '7.87'
$ws.Range('E31').Value = '  -1.94%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws 'D32' '13.68'
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('B33').Value = 'PEPE'
$ws.Range('C33').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws 'D33' '0.0000105'
$ws.Range('E33').Value = '  +17.43%  '
Set-TextValue $ws 'D34' '0.129'
$ws.Range('E34').Value = '  -3.84%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D35' '47.91'
$ws.Range('E35').Value = '  -2.82%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws 'D36' '70.99'
$ws.Range('E36').Value = '  +6.52%  '
Set-TextValue $ws 'D37' '646.11'
$ws.Range('E37').Value = '  -5.29%  '
Set-TextValue $ws 'D38' '0.429'
$ws.Range('E38').Value = '  -4.67%  '
Set-TextValue $ws 'D39' '3.42'
$ws.Range('E39').Value = '  +1.50%  '
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('E41').Value = '  -1.32%  '
$ws.Range('E42').Value = '  +0.09%  '
Set-TextValue $ws 'D43' '10.74'
$ws.Range('E43').Value = '  -3.88%  '
Set-TextValue $ws 'D44' '0.0482'
$ws.Range('E44').Value = '  -2.25%  '
Set-TextValue $ws 'D45' '3.16'
$ws.Range('E45').Value = '  -8.46%  '
$ws.Range('E46').Value = '  -1.86%  '
Set-TextValue $ws 'D47' '2.99'
$ws.Range('E47').Value = '  +31.28%  '
Set-TextValue $ws 'D48' '3.43'
$ws.Range('E48').Value = '  +3.27%  '
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue $ws 'D49' '0.000288'
$ws.Range('E49').Value = '  +2.65%  '
$ws.Range('B50').Value = 'Fetch.AI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws 'D50' '2.56'
$ws.Range('E50').Value = '  -4.03%  '
$ws.Range('D51').Value = '2.837.32'
$ws.Range('E51').Value = '  +3.54%  '

Write-Host "Updated cryptos list"
